$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cell: "PM-Sheet" -> "PM-Sheet (2)" ---
$ws.Range("B1").Value = "PM-Sheet (2)"

# --- Deadline date: 14-Mar-2019 -> 14-Feb-2019 ---
$ws.Range("D4").Value = 43510

# --- Row 12: "Add a registration window in the login area" task, now 100% done ---
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Add a registration window in the login area"
$ws.Range("D12").Value = "Djukic, Hamzic, Taha"
$ws.Range("E12").Value = "1 hour"
$ws.Range("F12").Value = "done"
$ws.Range("G12").Value = "1 hour"
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = "Djukic, Hamzic, Taha"

# --- Row 13: "Work on the marketing WebSite" task, now 100% done ---
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Work on the marketing WebSite"
$ws.Range("D13").Value = "Djukic, Hamzic, Taha"
$ws.Range("E13").Value = "6 hours"
$ws.Range("F13").Value = "done"
$ws.Range("G13").Value = "6 hours"
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = "Djukic, Hamzic, Taha"

# --- Row 14: clear out the old "android app" task entirely (stays 0%) ---
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("G14").ClearContents()
$ws.Range("H14").ClearContents()

# --- Row 15: clear out the old "marketing website" task entirely (stays 0%) ---
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("F15").Clear()
$ws.Range("G15").Clear()
$ws.Range("H15").Clear()

# --- Restore the current selection/view ---
$ws.Range("E12").Select()
